{"js": "// Replace the body text of 7 \"phishing message\" paragraphs with new\n// message content (per the diff). Each target paragraph consists of a\n// single run whose text is split into several sentences/lines joined by\n// manual line breaks (\\v, i.e. Word's <w:br/>). Office.js exposes a\n// manual line break as \"\\v\" (U+000B) inside Range/Paragraph text, so we\n// can replace a whole paragraph's content (text + breaks) in one\n// `insertText(..., \"Replace\")` call.\nconst edits = [\n  { index: 7, text: \"   Dear Esther Lindsley,\\u000b    We have detected a problem with your account. Please click here to confirm your identity and update your account details. If you did not request this email, please ignore it. Thank you for your cooperation.\\u000b    The Finance and Investment Department of the Bank of Romania\" },\n  { index: 10, text: \"Dear Mr. Sanders\\u000b\\u000bWe have your packing in our warehouse and we would like to ship but there's still outstanding amount of $10 for delivery please click on the link for payment. Thank you\" },\n  { index: 13, text: \"Amazon Prime are offering live-streaming of a wide range of sports from all over the world for a fraction of the normal price.\\u000b\\u000bSimply click the link below to find out more info!\" },\n  { index: 17, text: \"Dear Ms.Welling, \\u000b\\u000bThis message is to inform your hat your credit card with Novo Banco has been placed on hold due to suspicious charges. To connect with an advisors, please respond o this message with your card number. \\u000b\\u000bThank you,\\u000b\\u000bNovo Banco - Braganca\" },\n  { index: 19, text: \"Dear Charles Welling,\\u000b\\u000bWe recently noticed some unusual activity on your Amazon account. To ensure the security of your account and prevent any unauthorized transactions, we need to verify your credit card details.\\u000b\\u000bPlease reply to this message with the following information:\\u000b1. Your full name\\u000b2. Your 16-digit credit card number\\u000b3. The expiration date of your credit card\\u000b4. The 3-digit security code on the back of your card\\u000b\\u000bYour prompt response will help us secure your account and continue providing you with the best shopping experience.\\u000b\\u000bThank you for your cooperation.\\u000b\\u000bBest regards,\\u000bAmazon Security Team\" },\n  { index: 24, text: \"Hi Kellie, we are contacting you from the public library in Novokuybysjevsk. We see you have a few books you haven\u2019t returned to the library. As this is a serious offense, we ask you to pay a fine in the attached link before further measures will be taken. \" },\n  { index: 26, text: \"    Hello Kellie,\\u000b    We noticed some suspicious activity on your account. To keep your account secure, we need to verify your information. Please reply with the following information:\\u000b    - Credit card number: \\u000b    - Credit card expiration date: \\u000b    - Credit card CVV: \\u000b    If you did not initiate this request, please ignore this message.\\u000b    Thank you for your cooperation.\\u000b    Best regards,\\u000b    Novokuybyshevsk Bank\" },\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const edit of edits) {\n  const paragraph = paragraphs.items[edit.index];\n  paragraph.insertText(edit.text, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$d.Paragraphs.Item(8).Range.Text = '   Dear Esther Lindsley,' + \"`v\" + '    We have detected a problem with your account. Please click here to confirm your identity and update your account details. If you did not request this email, please ignore it. Thank you for your cooperation.' + \"`v\" + '    The Finance and Investment Department of the Bank of Romania'\n\n$d.Paragraphs.Item(11).Range.Text = 'Dear Mr. Sanders' + \"`v`v\" + 'We have your packing in our warehouse and we would like to ship but there''s still outstanding amount of $10 for delivery please click on the link for payment. Thank you'\n\n$d.Paragraphs.Item(14).Range.Text = 'Amazon Prime are offering live-streaming of a wide range of sports from all over the world for a fraction of the normal price.' + \"`v`v\" + 'Simply click the link below to find out more info!'\n\n$d.Paragraphs.Item(18).Range.Text = 'Dear Ms.Welling, ' + \"`v`v\" + 'This message is to inform your hat your credit card with Novo Banco has been placed on hold due to suspicious charges. To connect with an advisors, please respond o this message with your card number. ' + \"`v`v\" + 'Thank you,' + \"`v`v\" + 'Novo Banco - Braganca'\n\n$d.Paragraphs.Item(20).Range.Text = 'Dear Charles Welling,' + \"`v`v\" + 'We recently noticed some unusual activity on your Amazon account. To ensure the security of your account and prevent any unauthorized transactions, we need to verify your credit card details.' + \"`v`v\" + 'Please reply to this message with the following information:' + \"`v\" + '1. Your full name' + \"`v\" + '2. Your 16-digit credit card number' + \"`v\" + '3. The expiration date of your credit card' + \"`v\" + '4. The 3-digit security code on the back of your card' + \"`v`v\" + 'Your prompt response will help us secure your account and continue providing you with the best shopping experience.' + \"`v`v\" + 'Thank you for your cooperation.' + \"`v`v\" + 'Best regards,' + \"`v\" + 'Amazon Security Team'\n\n$d.Paragraphs.Item(25).Range.Text = 'Hi Kellie, we are contacting you from the public library in Novokuybysjevsk. We see you have a few books you haven\u2019t returned to the library. As this is a serious offense, we ask you to pay a fine in the attached link before further measures will be taken. '\n\n$d.Paragraphs.Item(27).Range.Text = '    Hello Kellie,' + \"`v\" + '    We noticed some suspicious activity on your account. To keep your account secure, we need to verify your information. Please reply with the following information:' + \"`v\" + '    - Credit card number: ' + \"`v\" + '    - Credit card expiration date: ' + \"`v\" + '    - Credit card CVV: ' + \"`v\" + '    If you did not initiate this request, please ignore this message.' + \"`v\" + '    Thank you for your cooperation.' + \"`v\" + '    Best regards,' + \"`v\" + '    Novokuybyshevsk Bank'\n"}
